$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Heat Exchangers": insert a new "Duty (MJ/h)" column after column B.
# Existing columns C:J (rows 1-7) shift right to D:K. Row 18 (C18, a lone
# styled/empty cell used for a border) is intentionally left untouched, so
# a plain Columns.Insert() (which would shift the whole column, row 18
# included) is avoided in favor of writing each target cell explicitly.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Heat Exchangers")

# --- Row 1 (headers) : shift C1:J1 -> D1:K1, then set new header C1 -------
$ws.Range("D1").Value = "Tin (C)"
$ws.Range("E1").Value = "Tout (C)"
$ws.Range("F1").Value = "Tutil (C)"
$ws.Range("G1").Value = "Tlm (C)"
$ws.Range("H1").Value = "util"
$ws.Range("I1").Value = "C ($/GJ)"
$ws.Range("J1").Value = "U (W/m2K)"
$ws.Range("K1").Value = "A"
$ws.Range("C1").Value = "Duty (MJ/h)"

# --- Row 2 ------------------------------------------------------------------
$ws.Range("D2").Value = 134.2
$ws.Range("E2").Value = 154.7
$ws.Range("F2").Value = 184
$ws.Range("G2").Formula = "=((F2-D2)-(F2-E2))/LN((F2-D2)/(F2-E2))"
$ws.Range("H2").Value = "mps"
$ws.Range("I2").Value = 14.19
$ws.Range("J2").Formula = "=400"
$ws.Range("K2").Formula = "=B2*1000/(G2*J2)"
$ws.Range("C2").Formula = "=B2*3.6"

# --- Row 3 ------------------------------------------------------------------
$ws.Range("D3").Value = 400
$ws.Range("E3").Value = 110
$ws.Range("F3").Value = 35
$ws.Range("G3").Formula = "=((F3-D3)-(F3-E3))/LN((F3-D3)/(F3-E3))"
$ws.Range("H3").Value = "cw"
$ws.Range("I3").Value = 0.354
$ws.Range("J3").Formula = "=400"
$ws.Range("K3").Formula = "=B3*1000/(G3*J3)"
$ws.Range("C3").Formula = "=B3*3.6"

# --- Row 4 ------------------------------------------------------------------
$ws.Range("D4").Value = 110
$ws.Range("E4").Value = 180
$ws.Range("F4").Value = 184
$ws.Range("G4").Formula = "=((F4-D4)-(F4-E4))/LN((F4-D4)/(F4-E4))"
$ws.Range("H4").Value = "mps"
$ws.Range("I4").Value = 14.19
$ws.Range("J4").Formula = "=400"
$ws.Range("K4").Formula = "=B4*1000/(G4*J4)"
$ws.Range("C4").Formula = "=B4*3.6"

# --- Row 5 ------------------------------------------------------------------
$ws.Range("D5").Value = 264
$ws.Range("E5").Value = 210
$ws.Range("F5").Value = 115
$ws.Range("G5").Formula = "=((F5-D5)-(F5-E5))/LN((F5-D5)/(F5-E5))"
$ws.Range("H5").Value = "bfw"
$ws.Range("I5").Value = 1.08
$ws.Range("J5").Formula = "=400"
$ws.Range("K5").Formula = "=B5*1000/(G5*J5)"
$ws.Range("C5").Formula = "=B5*3.6"

# --- Row 6 (note: original column C was blank, so new D is blank too) ------
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 159
$ws.Range("F6").Value = 115
$ws.Range("G6").Formula = "=F6-E6"
$ws.Range("H6").Value = "bfw"
$ws.Range("I6").Value = 1.08
$ws.Range("J6").Formula = "=400"
$ws.Range("K6").Formula = "=B6*1000/(G6*J6)"
$ws.Range("C6").Formula = "=B6*3.6"

# --- Row 7 (note: original column C was blank, so new D is blank too) ------
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = 215.3
$ws.Range("F7").Value = 115
$ws.Range("G7").Formula = "=F7-E7"
$ws.Range("H7").Value = "bfw"
$ws.Range("I7").Value = 1.08
$ws.Range("J7").Formula = "=400"
$ws.Range("K7").Formula = "=B7*1000/(G7*J7)"
$ws.Range("C7").Formula = "=B7*3.6"

# New column width (matches source layout) and selection.
$ws.Columns("C").ColumnWidth = 12.18
$ws.Range("B13:B15").Select()

# ---------------------------------------------------------------------------
# Sheet "Towers": update raw material flow for T-701 (B2): 33 -> 30.
# ---------------------------------------------------------------------------
$wsTowers = $wb.Worksheets.Item("Towers")
$wsTowers.Range("B2").Value = 30
$wsTowers.Range("E7").Select()
